$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert three new rows above the existing table data (rows 2-5 shift down
# to rows 5-8).
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).Insert()

# Match the formatting Excel already uses for the data rows (date style in
# column A, bordered body style in columns B/C) by copying it down from the
# first untouched data row.
$ws.Range("A5:C5").Copy()
$ws.Range("A2:C4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the three new "Sunday Worship" rows: date 9/28/2025, no song, topic
# "Sunday Worship".
$newDate = Get-Date -Year 2025 -Month 9 -Day 28 -Hour 0 -Minute 0 -Second 0
foreach ($r in 2..4) {
    $ws.Cells.Item($r, 1).Value = $newDate
    $ws.Cells.Item($r, 2).Value = ""
    $ws.Cells.Item($r, 3).Value = "Sunday Worship"
}

# Grow the table (and its autofilter) to cover the new rows.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:C8"))

$ws.Range("B2").Select()
